# Auto-generated script applying market-price/profit value updates
# to the Rafflesia_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# For each leve row, the currentAveragePrice / NQ / HQ / LevePrice / LeveProfit
# columns (H,I,J,K,L,M,N) are refreshed to newly observed market values.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 17
$ws.Range("H17").Value = 1692.5
$ws.Range("J17").Value = 1692.5
$ws.Range("L17").Value = 5077.5
$ws.Range("N17").Value = -5413.5

# Row 21
$ws.Range("H21").Value = 44000
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

# Row 23
$ws.Range("H23").Value = 44000
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

# Row 28
$ws.Range("H28").Value = 841
$ws.Range("I28").Value = 788.1667
$ws.Range("K28").Value = 788.1667
$ws.Range("M28").Value = -303.1667

# Row 40
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

# Row 43
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()

# Row 48
$ws.Range("H48").Value = 1019
$ws.Range("J48").Value = 1019
$ws.Range("L48").Value = 3057
$ws.Range("N48").Value = -3641

# Row 56
$ws.Range("H56").Value = 1019
$ws.Range("J56").Value = 1019
$ws.Range("L56").Value = 3057
$ws.Range("N56").Value = -4125

# Row 100
$ws.Range("H100").Value = 1995.5
$ws.Range("I100").Value = 1491
$ws.Range("J100").Value = 2500
$ws.Range("K100").Value = 1491
$ws.Range("L100").Value = 2500
$ws.Range("M100").Value = -950
$ws.Range("N100").Value = -3582

# Row 112
$ws.Range("H112").Value = 1382.6666
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1382.6666
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 4147.9998
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -6363.9998

# Row 137
$ws.Range("H137").Value = 3315.2727
$ws.Range("I137").Value = 2859.125
$ws.Range("K137").Value = 8577.375
$ws.Range("M137").Value = -6027.375

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 14
$ws.Range("H14").Value = 408.8
$ws.Range("J14").Value = 211
$ws.Range("L14").Value = 211
$ws.Range("N14").Value = -561

# Row 15
$ws.Range("H15").Value = 211
$ws.Range("J15").Value = 211
$ws.Range("L15").Value = 211
$ws.Range("N15").Value = -911

# Row 17
$ws.Range("H17").Value = 4000
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

# Row 61
$ws.Range("H61").Value = 11640.3125
$ws.Range("I61").Value = 7859.2856
$ws.Range("K61").Value = 7859.2856
$ws.Range("M61").Value = -7647.2856

# Row 74
$ws.Range("H74").Value = 4132.8335
$ws.Range("I74").Value = 3459.4
$ws.Range("K74").Value = 3459.4
$ws.Range("M74").Value = -2585.4

# Row 77
$ws.Range("H77").Value = 4132.8335
$ws.Range("I77").Value = 3459.4
$ws.Range("K77").Value = 17297
$ws.Range("M77").Value = -12929

# Row 132
$ws.Range("H132").Value = 2741.2856
$ws.Range("I132").Value = 1198.1666
$ws.Range("K132").Value = 3594.4998
$ws.Range("M132").Value = -1064.4998

# Row 136
$ws.Range("H136").Value = 11640.3125
$ws.Range("I136").Value = 7859.2856
$ws.Range("K136").Value = 23577.8568
$ws.Range("M136").Value = -21027.8568

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 14
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

# Row 86
$ws.Range("H86").Value = 3277.7778
$ws.Range("I86").Value = 1580
$ws.Range("K86").Value = 1580
$ws.Range("M86").Value = -457

# Row 89
$ws.Range("H89").Value = 3277.7778
$ws.Range("I89").Value = 1580
$ws.Range("K89").Value = 7900
$ws.Range("M89").Value = -2284

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 4
$ws.Range("H4").Value = 200
$ws.Range("I4").Value = 200
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 200
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -88
$ws.Range("N4").ClearContents()

# Row 6
$ws.Range("H6").Value = 5757875
$ws.Range("I6").Value = 6717104
$ws.Range("K6").Value = 6717104
$ws.Range("M6").Value = -6716991

# Row 122
$ws.Range("H122").Value = 3332.3333
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

# Row 132
$ws.Range("H132").Value = 4104.125
$ws.Range("J132").Value = 6194.5
$ws.Range("L132").Value = 18583.5
$ws.Range("N132").Value = -23643.5

# Row 134
$ws.Range("H134").Value = 2242.5715
$ws.Range("I134").Value = 2466.3333
$ws.Range("J134").Value = 900
$ws.Range("K134").Value = 7398.999899999999
$ws.Range("L134").Value = 2700
$ws.Range("M134").Value = -4863.999899999999
$ws.Range("N134").Value = -7770

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 12
$ws.Range("H12").Value = 26.5
$ws.Range("I12").Value = 19
$ws.Range("K12").Value = 57
$ws.Range("M12").Value = 116

# Row 23
$ws.Range("H23").Value = 556.8333
$ws.Range("J23").Value = 585
$ws.Range("L23").Value = 1755
$ws.Range("N23").Value = -2225

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 3
$ws.Range("H3").Value = 1412373.2
$ws.Range("I3").Value = 1336734.4
$ws.Range("J3").Value = 1444790
$ws.Range("K3").Value = 1336734.4
$ws.Range("L3").Value = 1444790
$ws.Range("M3").Value = -1336618.4
$ws.Range("N3").Value = -1445022

# Row 12
$ws.Range("H12").Value = 2750
$ws.Range("I12").Value = 5000
$ws.Range("J12").Value = 500
$ws.Range("K12").Value = 5000
$ws.Range("L12").Value = 500
$ws.Range("M12").Value = -4860
$ws.Range("N12").Value = -780

# Row 14
$ws.Range("H14").Value = 900
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 900
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 900
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -1236

# Row 80
$ws.Range("H80").Value = 36316.668
$ws.Range("I80").Value = 36000
$ws.Range("J80").Value = 36633.332
$ws.Range("K80").Value = 36000
$ws.Range("L80").Value = 36633.332
$ws.Range("M80").Value = -35002
$ws.Range("N80").Value = -38629.332

# Row 83
$ws.Range("H83").Value = 36316.668
$ws.Range("I83").Value = 36000
$ws.Range("J83").Value = 36633.332
$ws.Range("K83").Value = 180000
$ws.Range("L83").Value = 183166.66
$ws.Range("M83").Value = -175008
$ws.Range("N83").Value = -193150.66

# Row 94
$ws.Range("H94").Value = 18000
$ws.Range("J94").Value = 18000
$ws.Range("L94").Value = 18000
$ws.Range("N94").Value = -19352

# Row 132
$ws.Range("H132").Value = 2299.2222
$ws.Range("I132").Value = 1670.8572
$ws.Range("K132").Value = 5012.571599999999
$ws.Range("M132").Value = -2482.571599999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 4
$ws.Range("H4").Value = 1941.8
$ws.Range("I4").Value = 3654.5
$ws.Range("J4").Value = 800
$ws.Range("K4").Value = 3654.5
$ws.Range("L4").Value = 800
$ws.Range("M4").Value = -3541.5
$ws.Range("N4").Value = -1026

# Row 9
$ws.Range("H9").Value = 326.25
$ws.Range("I9").Value = 136
$ws.Range("J9").Value = 516.5
$ws.Range("K9").Value = 136
$ws.Range("L9").Value = 516.5
$ws.Range("M9").Value = 88
$ws.Range("N9").Value = -964.5

# Row 28
$ws.Range("H28").Value = 1941.8
$ws.Range("I28").Value = 3654.5
$ws.Range("J28").Value = 800
$ws.Range("K28").Value = 3654.5
$ws.Range("L28").Value = 800
$ws.Range("M28").Value = -3422.5
$ws.Range("N28").Value = -1264

# Row 37
$ws.Range("H37").Value = 1941.8
$ws.Range("I37").Value = 3654.5
$ws.Range("J37").Value = 800
$ws.Range("K37").Value = 3654.5
$ws.Range("L37").Value = 800
$ws.Range("M37").Value = -3547.5
$ws.Range("N37").Value = -1014

# Row 46
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

# Row 55
$ws.Range("H55").Value = 1132.9166
$ws.Range("I55").Value = 999
$ws.Range("K55").Value = 999
$ws.Range("M55").Value = -826

# Row 132
$ws.Range("H132").Value = 5048.8
$ws.Range("I132").Value = 4355.4287
$ws.Range("K132").Value = 13066.2861
$ws.Range("M132").Value = -10536.2861

# Row 136
$ws.Range("H136").Value = 29175.5
$ws.Range("I136").Value = 54452
$ws.Range("J136").Value = 3899
$ws.Range("K136").Value = 163356
$ws.Range("L136").Value = 11697
$ws.Range("M136").Value = -160806
$ws.Range("N136").Value = -16797

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 126
$ws.Range("H126").Value = 1000
$ws.Range("I126").Value = 1000
$ws.Range("K126").Value = 3000
$ws.Range("M126").Value = -530

# Row 132
$ws.Range("H132").Value = 1341.8334
$ws.Range("I132").Value = 1255.8889
$ws.Range("K132").Value = 3767.6667
$ws.Range("M132").Value = -1237.6667

# Row 135
$ws.Range("H135").Value = 41666.332
$ws.Range("J135").Value = 45000
$ws.Range("L135").Value = 45000
$ws.Range("N135").Value = -55140

# Row 136
$ws.Range("H136").Value = 1420.1578
$ws.Range("I136").Value = 1461.0834
$ws.Range("K136").Value = 4383.2502
$ws.Range("M136").Value = -1833.2502

